# Update the division-practice answer table: each table cell's answer
# text is replaced with its new value. Order matters for the last few
# replacements: the text produced by one Find/Replace must not be
# searched-for by a later one before it has had its turn (e.g. the
# "48÷9=5, 3" cell is renamed to "94÷9=10, 4" BEFORE a different cell is
# renamed to "48÷9=5, 3", otherwise the second rename would also hit the
# freshly-written text from the first).
$d = $word.ActiveDocument

$d.Content.Find.Execute("32÷4=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "18÷9=2, 0", 2) | Out-Null
$d.Content.Find.Execute("61÷5=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "55÷6=9, 1", 2) | Out-Null
$d.Content.Find.Execute("23÷7=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "55÷2=27, 1", 2) | Out-Null
$d.Content.Find.Execute("80÷5=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "18÷4=4, 2", 2) | Out-Null
$d.Content.Find.Execute("84÷4=21, 0", $true, $false, $false, $false, $false, $true, 1, $false, "26÷3=8, 2", 2) | Out-Null
$d.Content.Find.Execute("74÷3=24, 2", $true, $false, $false, $false, $false, $true, 1, $false, "94÷4=23, 2", 2) | Out-Null
$d.Content.Find.Execute("38÷9=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "24÷5=4, 4", 2) | Out-Null
$d.Content.Find.Execute("87÷8=10, 7", $true, $false, $false, $false, $false, $true, 1, $false, "50÷9=5, 5", 2) | Out-Null
$d.Content.Find.Execute("46÷6=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "58÷9=6, 4", 2) | Out-Null
$d.Content.Find.Execute("20÷5=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "22÷8=2, 6", 2) | Out-Null
$d.Content.Find.Execute("17÷9=1, 8", $true, $false, $false, $false, $false, $true, 1, $false, "10÷7=1, 3", 2) | Out-Null
$d.Content.Find.Execute("16÷4=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "32÷2=16, 0", 2) | Out-Null
$d.Content.Find.Execute("69÷4=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "70÷7=10, 0", 2) | Out-Null
$d.Content.Find.Execute("37÷2=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "40÷4=10, 0", 2) | Out-Null
$d.Content.Find.Execute("21÷4=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "50÷2=25, 0", 2) | Out-Null
$d.Content.Find.Execute("28÷5=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "98÷9=10, 8", 2) | Out-Null
$d.Content.Find.Execute("64÷5=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "21÷7=3, 0", 2) | Out-Null
$d.Content.Find.Execute("41÷6=6, 5", $true, $false, $false, $false, $false, $true, 1, $false, "71÷7=10, 1", 2) | Out-Null
$d.Content.Find.Execute("87÷7=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "94÷8=11, 6", 2) | Out-Null
$d.Content.Find.Execute("67÷3=22, 1", $true, $false, $false, $false, $false, $true, 1, $false, "66÷8=8, 2", 2) | Out-Null
$d.Content.Find.Execute("78÷9=8, 6", $true, $false, $false, $false, $false, $true, 1, $false, "95÷5=19, 0", 2) | Out-Null
$d.Content.Find.Execute("23÷3=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "40÷6=6, 4", 2) | Out-Null
$d.Content.Find.Execute("68÷2=34, 0", $true, $false, $false, $false, $false, $true, 1, $false, "55÷5=11, 0", 2) | Out-Null
$d.Content.Find.Execute("48÷9=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "94÷9=10, 4", 2) | Out-Null
$d.Content.Find.Execute("38÷6=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "48÷9=5, 3", 2) | Out-Null
